$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16
$ws.Cells.Item(16, 2).Value2 = 7100665
$ws.Cells.Item(16, 6).Value2 = "OFI Crete"
$ws.Cells.Item(16, 7).Value2 = "PAOK Salonika"
$ws.Cells.Item(16, 8).Value2 = 1
$ws.Cells.Item(16, 9).Value2 = 0
$ws.Cells.Item(16, 11).Value2 = 4.75
$ws.Cells.Item(16, 12).Value2 = 3.6
$ws.Cells.Item(16, 13).Value2 = 1.75
$ws.Cells.Item(16, 14).Value2 = 4.5
$ws.Cells.Item(16, 16).Value2 = 1.75
$ws.Cells.Item(16, 17).Value2 = 0.75
$ws.Cells.Item(16, 18).Value2 = 1.85
$ws.Cells.Item(16, 19).Value2 = 2
$ws.Cells.Item(16, 20).Value2 = 2.75
$ws.Cells.Item(16, 21).Value2 = 2
$ws.Cells.Item(16, 22).Value2 = 1.85
$ws.Cells.Item(16, 23).Value2 = 3.5
$ws.Cells.Item(16, 26).Value2 = 0.8500000000000001
$ws.Cells.Item(16, 27).Value2 = -1
$ws.Cells.Item(16, 28).Value2 = -1
$ws.Cells.Item(16, 29).Value2 = 0.8500000000000001

# Row 17
$ws.Cells.Item(17, 2).Value2 = 7100661
$ws.Cells.Item(17, 6).Value2 = "Aris Salonika"
$ws.Cells.Item(17, 7).Value2 = "Asteras Tripolis"
$ws.Cells.Item(17, 8).Value2 = 3
$ws.Cells.Item(17, 9).Value2 = 2
$ws.Cells.Item(17, 11).Value2 = 1.8
$ws.Cells.Item(17, 12).Value2 = 3.4
$ws.Cells.Item(17, 13).Value2 = 4.75
$ws.Cells.Item(17, 14).Value2 = 1.55
$ws.Cells.Item(17, 16).Value2 = 7
$ws.Cells.Item(17, 17).Value2 = -1
$ws.Cells.Item(17, 18).Value2 = 2
$ws.Cells.Item(17, 19).Value2 = 1.85
$ws.Cells.Item(17, 20).Value2 = 2.25
$ws.Cells.Item(17, 21).Value2 = 1.825
$ws.Cells.Item(17, 22).Value2 = 2.025
$ws.Cells.Item(17, 23).Value2 = 0.55
$ws.Cells.Item(17, 26).Value2 = 0
$ws.Cells.Item(17, 27).Value2 = 0
$ws.Cells.Item(17, 28).Value2 = 0.825
$ws.Cells.Item(17, 29).Value2 = -1

# Row 143
$ws.Cells.Item(143, 2).Value2 = 6937250
$ws.Cells.Item(143, 6).Value2 = "Giannina"
$ws.Cells.Item(143, 7).Value2 = "Lamia"
$ws.Cells.Item(143, 8).Value2 = 1
$ws.Cells.Item(143, 9).Value2 = 4
$ws.Cells.Item(143, 10).Value2 = "A"
$ws.Cells.Item(143, 11).Value2 = 2.3
$ws.Cells.Item(143, 12).Value2 = 3.25
$ws.Cells.Item(143, 13).Value2 = 3.25
$ws.Cells.Item(143, 14).Value2 = 2.55
$ws.Cells.Item(143, 15).Value2 = 2.875
$ws.Cells.Item(143, 16).Value2 = 3.1
$ws.Cells.Item(143, 17).Value2 = 0
$ws.Cells.Item(143, 18).Value2 = 1.75
$ws.Cells.Item(143, 19).Value2 = 2.125
$ws.Cells.Item(143, 20).Value2 = 2
$ws.Cells.Item(143, 21).Value2 = 1.85
$ws.Cells.Item(143, 22).Value2 = 2
$ws.Cells.Item(143, 23).Value2 = -1
$ws.Cells.Item(143, 25).Value2 = 2.1
$ws.Cells.Item(143, 26).Value2 = -1
$ws.Cells.Item(143, 27).Value2 = 1.125
$ws.Cells.Item(143, 28).Value2 = 0.8500000000000001

# Row 144
$ws.Cells.Item(144, 2).Value2 = 6937247
$ws.Cells.Item(144, 6).Value2 = "AEK Athens"
$ws.Cells.Item(144, 7).Value2 = "Asteras Tripolis"
$ws.Cells.Item(144, 8).Value2 = 4
$ws.Cells.Item(144, 9).Value2 = 2
$ws.Cells.Item(144, 10).Value2 = "H"
$ws.Cells.Item(144, 11).Value2 = 1.285
$ws.Cells.Item(144, 12).Value2 = 5.5
$ws.Cells.Item(144, 13).Value2 = 12
$ws.Cells.Item(144, 14).Value2 = 1.285
$ws.Cells.Item(144, 15).Value2 = 5.75
$ws.Cells.Item(144, 16).Value2 = 10
$ws.Cells.Item(144, 17).Value2 = -1.5
$ws.Cells.Item(144, 18).Value2 = 1.825
$ws.Cells.Item(144, 19).Value2 = 2.025
$ws.Cells.Item(144, 20).Value2 = 3
$ws.Cells.Item(144, 21).Value2 = 2.025
$ws.Cells.Item(144, 22).Value2 = 1.825
$ws.Cells.Item(144, 23).Value2 = 0.2849999999999999
$ws.Cells.Item(144, 25).Value2 = -1
$ws.Cells.Item(144, 26).Value2 = 0.825
$ws.Cells.Item(144, 27).Value2 = -1
$ws.Cells.Item(144, 28).Value2 = 1.025

# Row 168
$ws.Cells.Item(168, 2).Value2 = 6937267
$ws.Cells.Item(168, 6).Value2 = "Volos NFC"
$ws.Cells.Item(168, 7).Value2 = "OFI Crete"
$ws.Cells.Item(168, 9).Value2 = 1
$ws.Cells.Item(168, 10).Value2 = "H"
$ws.Cells.Item(168, 11).Value2 = 2.7
$ws.Cells.Item(168, 12).Value2 = 3.25
$ws.Cells.Item(168, 13).Value2 = 2.625
$ws.Cells.Item(168, 14).Value2 = 2.7
$ws.Cells.Item(168, 15).Value2 = 3.2
$ws.Cells.Item(168, 16).Value2 = 2.8
$ws.Cells.Item(168, 17).Value2 = 0
$ws.Cells.Item(168, 18).Value2 = 1.825
$ws.Cells.Item(168, 19).Value2 = 2.025
$ws.Cells.Item(168, 20).Value2 = 2.25
$ws.Cells.Item(168, 21).Value2 = 2
$ws.Cells.Item(168, 22).Value2 = 1.85
$ws.Cells.Item(168, 23).Value2 = 1.7
$ws.Cells.Item(168, 24).Value2 = -1
$ws.Cells.Item(168, 26).Value2 = 0.825
$ws.Cells.Item(168, 27).Value2 = -1
$ws.Cells.Item(168, 28).Value2 = 1

# Row 169
$ws.Cells.Item(169, 2).Value2 = 6935703
$ws.Cells.Item(169, 6).Value2 = "Asteras Tripolis"
$ws.Cells.Item(169, 7).Value2 = "Kifisias FC"
$ws.Cells.Item(169, 9).Value2 = 3
$ws.Cells.Item(169, 10).Value2 = "D"
$ws.Cells.Item(169, 11).Value2 = 1.833
$ws.Cells.Item(169, 12).Value2 = 3.4
$ws.Cells.Item(169, 13).Value2 = 4.5
$ws.Cells.Item(169, 14).Value2 = 1.8
$ws.Cells.Item(169, 15).Value2 = 3.5
$ws.Cells.Item(169, 16).Value2 = 4.75
$ws.Cells.Item(169, 17).Value2 = -0.75
$ws.Cells.Item(169, 18).Value2 = 2.025
$ws.Cells.Item(169, 19).Value2 = 1.825
$ws.Cells.Item(169, 20).Value2 = 2.5
$ws.Cells.Item(169, 21).Value2 = 1.875
$ws.Cells.Item(169, 22).Value2 = 1.975
$ws.Cells.Item(169, 23).Value2 = -1
$ws.Cells.Item(169, 24).Value2 = 2.5
$ws.Cells.Item(169, 26).Value2 = -1
$ws.Cells.Item(169, 27).Value2 = 0.825
$ws.Cells.Item(169, 28).Value2 = 0.875

# Row 176
$ws.Cells.Item(176, 2).Value2 = 6935700
$ws.Cells.Item(176, 6).Value2 = "Panserraikos"
$ws.Cells.Item(176, 7).Value2 = "Asteras Tripolis"
$ws.Cells.Item(176, 9).Value2 = 1
$ws.Cells.Item(176, 10).Value2 = "H"
$ws.Cells.Item(176, 11).Value2 = 2.6
$ws.Cells.Item(176, 12).Value2 = 3.2
$ws.Cells.Item(176, 13).Value2 = 2.875
$ws.Cells.Item(176, 14).Value2 = 2.25
$ws.Cells.Item(176, 16).Value2 = 3.3
$ws.Cells.Item(176, 17).Value2 = -0.25
$ws.Cells.Item(176, 18).Value2 = 1.925
$ws.Cells.Item(176, 19).Value2 = 1.925
$ws.Cells.Item(176, 21).Value2 = 2
$ws.Cells.Item(176, 22).Value2 = 1.85
$ws.Cells.Item(176, 23).Value2 = 1.25
$ws.Cells.Item(176, 24).Value2 = -1
$ws.Cells.Item(176, 26).Value2 = 0.925
$ws.Cells.Item(176, 27).Value2 = -1
$ws.Cells.Item(176, 28).Value2 = 1

# Row 177
$ws.Cells.Item(177, 2).Value2 = 6935701
$ws.Cells.Item(177, 6).Value2 = "Kifisias FC"
$ws.Cells.Item(177, 7).Value2 = "Panetolikos"
$ws.Cells.Item(177, 8).Value2 = 2
$ws.Cells.Item(177, 10).Value2 = "D"
$ws.Cells.Item(177, 11).Value2 = 2.45
$ws.Cells.Item(177, 12).Value2 = 3.25
$ws.Cells.Item(177, 13).Value2 = 3
$ws.Cells.Item(177, 14).Value2 = 2.05
$ws.Cells.Item(177, 15).Value2 = 3.3
$ws.Cells.Item(177, 16).Value2 = 3.8
$ws.Cells.Item(177, 17).Value2 = -0.5
$ws.Cells.Item(177, 18).Value2 = 2.05
$ws.Cells.Item(177, 19).Value2 = 1.8
$ws.Cells.Item(177, 20).Value2 = 2.25
$ws.Cells.Item(177, 21).Value2 = 1.8
$ws.Cells.Item(177, 22).Value2 = 2.05
$ws.Cells.Item(177, 24).Value2 = 2.3
$ws.Cells.Item(177, 25).Value2 = -1
$ws.Cells.Item(177, 27).Value2 = 0.8
$ws.Cells.Item(177, 28).Value2 = 0.8
$ws.Cells.Item(177, 29).Value2 = -1

# Row 178
$ws.Cells.Item(178, 2).Value2 = 6937272
$ws.Cells.Item(178, 6).Value2 = "Lamia"
$ws.Cells.Item(178, 7).Value2 = "PAOK Salonika"
$ws.Cells.Item(178, 8).Value2 = 0
$ws.Cells.Item(178, 9).Value2 = 2
$ws.Cells.Item(178, 10).Value2 = "A"
$ws.Cells.Item(178, 11).Value2 = 7.5
$ws.Cells.Item(178, 12).Value2 = 4.5
$ws.Cells.Item(178, 13).Value2 = 1.444
$ws.Cells.Item(178, 14).Value2 = 9.5
$ws.Cells.Item(178, 15).Value2 = 5
$ws.Cells.Item(178, 16).Value2 = 1.333
$ws.Cells.Item(178, 17).Value2 = 1.5
$ws.Cells.Item(178, 20).Value2 = 3
$ws.Cells.Item(178, 21).Value2 = 1.95
$ws.Cells.Item(178, 22).Value2 = 1.9
$ws.Cells.Item(178, 23).Value2 = -1
$ws.Cells.Item(178, 25).Value2 = 0.333
$ws.Cells.Item(178, 26).Value2 = -1
$ws.Cells.Item(178, 27).Value2 = 0.925
$ws.Cells.Item(178, 28).Value2 = -1
$ws.Cells.Item(178, 29).Value2 = 0.8999999999999999

# Row 212
$ws.Cells.Item(212, 2).Value2 = 7920482
$ws.Cells.Item(212, 5).Value2 = 45403.5625
$ws.Cells.Item(212, 6).Value2 = "Aris Salonika"
$ws.Cells.Item(212, 7).Value2 = "AEK Athens"
$ws.Cells.Item(212, 11).Value2 = 5.5
$ws.Cells.Item(212, 12).Value2 = 4.75
$ws.Cells.Item(212, 13).Value2 = 1.5
$ws.Cells.Item(212, 14).Value2 = 9.5
$ws.Cells.Item(212, 15).Value2 = 4.75
$ws.Cells.Item(212, 16).Value2 = 1.333
$ws.Cells.Item(212, 17).Value2 = 1.5
$ws.Cells.Item(212, 20).Value2 = 2.75
$ws.Cells.Item(212, 21).Value2 = 1.975
$ws.Cells.Item(212, 22).Value2 = 1.875

# Row 213
$ws.Cells.Item(213, 2).Value2 = 7920483
$ws.Cells.Item(213, 5).Value2 = 45403.60416666666
$ws.Cells.Item(213, 6).Value2 = "Olympiakos"
$ws.Cells.Item(213, 7).Value2 = "PAOK Salonika"
$ws.Cells.Item(213, 11).Value2 = 2.1
$ws.Cells.Item(213, 12).Value2 = 3.4
$ws.Cells.Item(213, 13).Value2 = 3.4
$ws.Cells.Item(213, 14).Value2 = 2.3
$ws.Cells.Item(213, 16).Value2 = 3
$ws.Cells.Item(213, 17).Value2 = -0.25
$ws.Cells.Item(213, 18).Value2 = 1.975
$ws.Cells.Item(213, 19).Value2 = 1.875
$ws.Cells.Item(213, 20).Value2 = 2.75
$ws.Cells.Item(213, 21).Value2 = 1.85
$ws.Cells.Item(213, 22).Value2 = 2

# Remove obsolete trailing fixtures (rows 214-218), data refreshed/reordered upstream
$ws.Rows("214:218").Delete()
